$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1300
$ws.Range("I2").Value = 1067
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 1067
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -954
$ws.Range("N2").Value = -2225
# Row 18
$ws.Range("H18").Value = 790.8
$ws.Range("I18").Value = 790.8
$ws.Range("K18").Value = 790.8
$ws.Range("M18").Value = -506.8
# Row 51
$ws.Range("H51").Value = 35062.125
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 38642.43
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 38642.43
$ws.Range("M51").Value = -9516
$ws.Range("N51").Value = -39610.43
# Row 116
$ws.Range("H116").Value = 6856.6665
$ws.Range("J116").Value = 7164.5557
$ws.Range("L116").Value = 7164.5557
$ws.Range("N116").Value = -14048.5557
# Row 127
$ws.Range("H127").Value = 705
$ws.Range("I127").Value = 705
$ws.Range("K127").Value = 2115
$ws.Range("M127").Value = 2845
# Row 132
$ws.Range("H132").Value = 2747.0278
$ws.Range("I132").Value = 2810.6775
$ws.Range("J132").Value = 2352.4
$ws.Range("K132").Value = 8432.032499999999
$ws.Range("L132").Value = 7057.200000000001
$ws.Range("M132").Value = -5902.032499999999
$ws.Range("N132").Value = -12117.2
# Row 137
$ws.Range("H137").Value = 2066.4
$ws.Range("I137").Value = 1167.1
$ws.Range("J137").Value = 2965.7
$ws.Range("K137").Value = 3501.3
$ws.Range("L137").Value = 8897.099999999999
$ws.Range("M137").Value = -951.2999999999997
$ws.Range("N137").Value = -13997.1
# Row 138
$ws.Range("H138").Value = 3639.5083
$ws.Range("I138").Value = 1305.6522
$ws.Range("J138").Value = 5052.1055
$ws.Range("K138").Value = 3916.9566
$ws.Range("L138").Value = 15156.3165
$ws.Range("M138").Value = 1223.0434
$ws.Range("N138").Value = -25436.3165
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 29500
$ws.Range("I28").Value = 29500
$ws.Range("K28").Value = 29500
$ws.Range("M28").Value = -29308
# Row 61
$ws.Range("H61").Value = 7903.4546
$ws.Range("I61").Value = 7144.25
$ws.Range("J61").Value = 9928
$ws.Range("K61").Value = 7144.25
$ws.Range("L61").Value = 9928
$ws.Range("M61").Value = -6932.25
$ws.Range("N61").Value = -10352
# Row 99
$ws.Range("H99").Value = 29500
$ws.Range("I99").Value = 29500
$ws.Range("K99").Value = 29500
$ws.Range("M99").Value = -26505
# Row 134
$ws.Range("H134").Value = 69999.5
$ws.Range("J134").Value = 69999.5
$ws.Range("L134").Value = 69999.5
$ws.Range("N134").Value = -80139.5
# Row 135
$ws.Range("H135").Value = 60249.668
$ws.Range("J135").Value = 60249.668
$ws.Range("L135").Value = 60249.668
$ws.Range("N135").Value = -70389.66800000001
# Row 136
$ws.Range("H136").Value = 7903.4546
$ws.Range("I136").Value = 7144.25
$ws.Range("J136").Value = 9928
$ws.Range("K136").Value = 21432.75
$ws.Range("L136").Value = 29784
$ws.Range("M136").Value = -18882.75
$ws.Range("N136").Value = -34884

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 55
$ws.Range("H55").Value = 400000
$ws.Range("J55").Value = 400000
$ws.Range("L55").Value = 400000
$ws.Range("N55").Value = -400546
# Row 86
$ws.Range("H86").Value = 3508.6
$ws.Range("I86").Value = 3233.5908
$ws.Range("J86").Value = 4264.875
$ws.Range("K86").Value = 3233.5908
$ws.Range("L86").Value = 4264.875
$ws.Range("M86").Value = -2110.5908
$ws.Range("N86").Value = -6510.875
# Row 89
$ws.Range("H89").Value = 3508.6
$ws.Range("I89").Value = 3233.5908
$ws.Range("J89").Value = 4264.875
$ws.Range("K89").Value = 16167.954
$ws.Range("L89").Value = 21324.375
$ws.Range("M89").Value = -10551.954
$ws.Range("N89").Value = -32556.375
# Row 94
$ws.Range("H94").Value = 890.5714
$ws.Range("J94").Value = 1375
$ws.Range("L94").Value = 1375
$ws.Range("N94").Value = -2277
# Row 134
$ws.Range("H134").Value = 3591.1333
$ws.Range("I134").Value = 3591.1333
$ws.Range("K134").Value = 10773.3999
$ws.Range("M134").Value = -8238.3999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 3598.9048
$ws.Range("I132").Value = 2987.9443
$ws.Range("K132").Value = 8963.832900000001
$ws.Range("M132").Value = -6433.832900000001
# Row 134
$ws.Range("H134").Value = 5925.4243
$ws.Range("I134").Value = 5541.483
$ws.Range("J134").Value = 8709
$ws.Range("K134").Value = 16624.449
$ws.Range("L134").Value = 26127
$ws.Range("M134").Value = -14089.449
$ws.Range("N134").Value = -31197

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 20
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -3454
# Row 21
$ws.Range("H21").Value = 5616.3335
$ws.Range("I21").Value = 2925
$ws.Range("J21").Value = 10999
$ws.Range("K21").Value = 8775
$ws.Range("L21").Value = 32997
$ws.Range("M21").Value = -8602
$ws.Range("N21").Value = -33343
# Row 23
$ws.Range("H23").Value = 774.913
$ws.Range("J23").Value = 830
$ws.Range("L23").Value = 2490
$ws.Range("N23").Value = -2960
# Row 47
$ws.Range("H47").Value = 2444
$ws.Range("I47").Value = 733
$ws.Range("K47").Value = 2199
$ws.Range("M47").Value = -1768
# Row 75
$ws.Range("H75").Value = 2240.4546
$ws.Range("J75").Value = 2173.1428
$ws.Range("L75").Value = 6519.428400000001
$ws.Range("N75").Value = -8515.428400000001
# Row 76
$ws.Range("H76").Value = 17940
$ws.Range("I76").Value = 8526.666999999999
$ws.Range("K76").Value = 25580.001
$ws.Range("M76").Value = -25197.001
# Row 78
$ws.Range("H78").Value = 2240.4546
$ws.Range("J78").Value = 2173.1428
$ws.Range("L78").Value = 19558.2852
$ws.Range("N78").Value = -29542.2852
# Row 79
$ws.Range("H79").Value = 17940
$ws.Range("I79").Value = 8526.666999999999
$ws.Range("K79").Value = 25580.001
$ws.Range("M79").Value = -24254.001
# Row 126
$ws.Range("H126").Value = 2369.8
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 1644.3636
$ws.Range("I132").Value = 677.8
$ws.Range("K132").Value = 6100.2
$ws.Range("M132").Value = -3570.2
# Row 133
$ws.Range("H133").Value = 11570.143
$ws.Range("I133").Value = 6198.2
$ws.Range("K133").Value = 18594.6
$ws.Range("M133").Value = -13534.6

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 3422499.8
$ws.Range("I3").Value = 6674333.5
$ws.Range("J3").Value = 170666.33
$ws.Range("K3").Value = 6674333.5
$ws.Range("L3").Value = 170666.33
$ws.Range("M3").Value = -6674217.5
$ws.Range("N3").Value = -170898.33
# Row 10
$ws.Range("H10").Value = 2009580.2
$ws.Range("I10").Value = 5002001
$ws.Range("J10").Value = 14633
$ws.Range("K10").Value = 5002001
$ws.Range("L10").Value = 14633
$ws.Range("M10").Value = -5001832
$ws.Range("N10").Value = -14971
# Row 14
$ws.Range("H14").Value = 130334.25
$ws.Range("I14").Value = 167362.5
$ws.Range("K14").Value = 167362.5
$ws.Range("M14").Value = -167194.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 19
$ws.Range("H19").Value = 10002
$ws.Range("I19").Value = 10002
$ws.Range("K19").Value = 10002
$ws.Range("M19").Value = -9832
# Row 122
$ws.Range("H122").Value = 12687.8
$ws.Range("I122").Value = 12687.8
$ws.Range("K122").Value = 38063.39999999999
$ws.Range("M122").Value = -35613.39999999999
# Row 141
$ws.Range("H141").Value = 84147.30499999999
$ws.Range("J141").Value = 84147.30499999999
$ws.Range("L141").Value = 84147.30499999999
$ws.Range("N141").Value = -94507.30499999999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 2150
$ws.Range("J7").Value = 3600
$ws.Range("L7").Value = 3600
$ws.Range("N7").Value = -3826
# Row 12
$ws.Range("H12").Value = 7551
$ws.Range("I12").Value = 7551
$ws.Range("K12").Value = 7551
$ws.Range("M12").Value = -7409
# Row 107
$ws.Range("H107").Value = 4038.2
$ws.Range("I107").Value = 5386.25
$ws.Range("J107").Value = 2497.5715
$ws.Range("K107").Value = 16158.75
$ws.Range("L107").Value = 7492.7145
$ws.Range("M107").Value = -14238.75
$ws.Range("N107").Value = -11332.7145
# Row 132
$ws.Range("H132").Value = 1870.5193
$ws.Range("I132").Value = 1852.4359
$ws.Range("K132").Value = 5557.307699999999
$ws.Range("M132").Value = -3027.307699999999
